$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("D-P")
$ws.Range("E17").Value = 11
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = 8
$ws.Range("G19").Value = 10
$ws.Activate()
$ws.Range("K17").Select()
